$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$val) {
    if ($val -match '^[+-]?\d+(\.\d+)?$') {
        $cell.Value = "'" + $val
    } else {
        $cell.Value = $val
    }
}

$changes = @(
    @{Row=2; D="33.590.37"; E="  +8.85%  "},
    @{Row=3; D="1.782.92"; E="  +6.24%  "},
    @{Row=4; E="  +0.99%  "},
    @{Row=5; D="224.53"; E="  +2.30%  "},
    @{Row=6; D="0.555"; E="  +3.64%  "},
    @{Row=7; E="  +0.80%  "},
    @{Row=8; D="31.02"; E="  +7.02%  "},
    @{Row=9; D="46.23"; E="  +4.75%  "},
    @{Row=10; D="0.279"; E="  +5.71%  "},
    @{Row=11; D="0.0654"; E="  +1.87%  "},
    @{Row=12; D="0.0927"; E="  +2.40%  "},
    @{Row=13; D="2.055.19"; E="  +6.96%  "},
    @{Row=14; D="1.805.11"; E="  +7.39%  "},
    @{Row=15; D="0.628"; E="  +4.21%  "},
    @{Row=16; D="33.743.02"; E="  +9.39%  "},
    @{Row=17; D="9.89"; E="  -1.76%  "},
    @{Row=18; D="4.17"; E="  +3.02%  "},
    @{Row=19; D="68.42"; E="  +3.84%  "},
    @{Row=20; D="249.02"; E="  +2.50%  "},
    @{Row=21; D="0.0₃0734"; E="  +2.15%  "},
    @{Row=22; E="  +0.48%  "},
    @{Row=23; D="10.29"; E="  +3.25%  "},
    @{Row=24; D="4.21"; E="  -0.37%  "},
    @{Row=25; D="2.16"; E="  -0.14%  "},
    @{Row=26; D="158.06"; E="  -0.55%  "},
    @{Row=27; D="16.30"; E="  +3.06%  "},
    @{Row=28; D="0.114"; E="  +1.66%  "},
    @{Row=29; D="6.87"; E="  +3.03%  "},
    @{Row=30; D="1.00"; E="  +0.60%  "},
    @{Row=31; D="3.77"; E="  +9.09%  "},
    @{Row=32; D="0.0509"; E="  +3.46%  "},
    @{Row=33; D="1.19"; E="  +4.53%  "},
    @{Row=34; D="3.50"; E="  +6.08%  "},
    @{Row=35; D="1.491.14"; E="  -1.55%  "},
    @{Row=36; D="1.73"; E="  -0.81%  "},
    @{Row=37; D="1.06"; E="  +4.47%  "},
    @{Row=38; D="0.0185"; E="  +3.71%  "},
    @{Row=39; D="0.617"; E="  +1.66%  "},
    @{Row=40; D="82.04"; E="  -2.24%  "},
    @{Row=41; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="2.77"; E="  +4.66%  "},
    @{Row=42; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="2.38"; E="  +3.72%  "},
    @{Row=43; D="0.881"; E="  +5.08%  "},
    @{Row=44; D="2.06"; E="  +1.18%  "},
    @{Row=45; D="0.0513"; E="  +2.75%  "},
    @{Row=46; D="1.07"; E="  +4.50%  "},
    @{Row=47; E="  +7.10%  "},
    @{Row=48; B="PaxDollar"; C="https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; D="1.01"; E="  +0.59%  "},
    @{Row=49; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="5.71"; E="  +2.64%  "},
    @{Row=50; D="11.68"; E="  +11.20%  "},
    @{Row=51; D="50.28"; E="  -0.75%  "}
)

foreach ($item in $changes) {
    if ($item.ContainsKey('B')) { $ws.Cells.Item($item.Row, 2).Value = $item.B }
    if ($item.ContainsKey('C')) { $ws.Cells.Item($item.Row, 3).Value = $item.C }
    if ($item.ContainsKey('D')) { Set-TextValue $ws.Cells.Item($item.Row, 4) $item.D }
    if ($item.ContainsKey('E')) { $ws.Cells.Item($item.Row, 5).Value = $item.E }
}
